$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a brand-new first paragraph containing "John Brandenburg"
#    and carry the "_GoBack" bookmark onto the end of it.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(1)
# Add a trailing placeholder character so the bookmark insertion point
# below is not exactly on a paragraph/text boundary (a boundary offset
# gets mis-anchored to the very start of the document by this host).
$newPara.Range.Text = "John BrandenburgX"

$bmPos = $d.Range(16, 16)
$d.Bookmarks.Add("_GoBack", $bmPos)

# Remove the placeholder character now that the bookmark is anchored.
$placeholder = $d.Range(16, 17)
$placeholder.Delete()

# ------------------------------------------------------------------
# 2) Collapse the hyperlink's two runs (and the bookmark that used to
#    sit between them) into a single run reading "Project 2 GH link".
#    Changing TextToDisplay rewrites the hyperlink's run cleanly.
# ------------------------------------------------------------------
$h = $d.Hyperlinks.Item(1)
$h.TextToDisplay = "Project 2 GH link"
